$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 80
$ws1.Range("F3").Value = 3961
$ws1.Range("F4").Value = 2327
$ws1.Range("F5").Value = 463
$ws1.Range("F7").Value = 29
$ws1.Range("F8").Value = 10
$ws1.Range("F10").Value = 116
$ws1.Range("F11").Value = 35
$ws1.Range("F12").Value = 121
$ws1.Range("F13").Value = 1469
$ws1.Range("F14").Value = 260
$ws1.Range("F15").Value = 2704
$ws1.Range("F16").Value = 184

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 37

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 80
$ws4.Range("F3").Value = 3961
$ws4.Range("F4").Value = 2327
$ws4.Range("F5").Value = 463
$ws4.Range("F7").Value = 29
$ws4.Range("F8").Value = 10
$ws4.Range("F9").Value = 37
$ws4.Range("F11").Value = 116
$ws4.Range("F12").Value = 35
$ws4.Range("F13").Value = 121
$ws4.Range("F16").Value = 1469
$ws4.Range("F17").Value = 260
$ws4.Range("F18").Value = 2704
$ws4.Range("F19").Value = 184
